$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 61; $r++) {
    $bCell = $ws.Cells.Item($r, 2)   # column B (Sales)
    $eCell = $ws.Cells.Item($r, 5)   # column E (Ad_Spend_X)
    $oldB = $bCell.Value2
    $eVal = $eCell.Value2
    $newB = $oldB - (0.8 * $eVal)
    $bCell.Value = $newB
}
